$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting (date number format/style) from the last existing
# data row (A18) into the new row, then set the actual values.
$ws.Range("A18").Copy($ws.Range("A19"))

$ws.Range("A19").Value = 45897
$ws.Range("B19").Value = 4
